$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 300 (shifts existing rows 300+ down by 3),
# matching the target diff where a new block of 3 price rows (date 44468)
# is inserted and everything below shifts down.
$ws.Range("A300:A302").EntireRow.Insert()

# Row 300: new "1a/2a/3a amarillo" price entry for fecha 44468
$ws.Cells.Item(300,1).Value = 2
$ws.Cells.Item(300,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(300,3).Value = "Coquimbo"
$ws.Cells.Item(300,4).Value = 44468
$ws.Cells.Item(300,5).Value = 4
$ws.Cells.Item(300,6).Value = "Fruta"
$ws.Cells.Item(300,7).Value = 100102
$ws.Cells.Item(300,8).Value = "Cítricos"
$ws.Cells.Item(300,9).Value = 100102003
$ws.Cells.Item(300,10).Value = "Limón"
$ws.Cells.Item(300,11).Value = "Sin especificar"
$ws.Cells.Item(300,12).Value = "1a amarillo"
$ws.Cells.Item(300,13).Value = 750
$ws.Cells.Item(300,14).Value = 4300
$ws.Cells.Item(300,15).Value = 4500
$ws.Cells.Item(300,16).Value = 4400
$ws.Cells.Item(300,17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(300,18).Value = "Provincia de Limarí"
$ws.Cells.Item(300,19).Value = 275
$ws.Cells.Item(300,20).Value = 16

# Row 301: new "1a/2a/3a amarillo" price entry for fecha 44468
$ws.Cells.Item(301,1).Value = 2
$ws.Cells.Item(301,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(301,3).Value = "Coquimbo"
$ws.Cells.Item(301,4).Value = 44468
$ws.Cells.Item(301,5).Value = 4
$ws.Cells.Item(301,6).Value = "Fruta"
$ws.Cells.Item(301,7).Value = 100102
$ws.Cells.Item(301,8).Value = "Cítricos"
$ws.Cells.Item(301,9).Value = 100102003
$ws.Cells.Item(301,10).Value = "Limón"
$ws.Cells.Item(301,11).Value = "Sin especificar"
$ws.Cells.Item(301,12).Value = "2a amarillo"
$ws.Cells.Item(301,13).Value = 600
$ws.Cells.Item(301,14).Value = 3300
$ws.Cells.Item(301,15).Value = 3500
$ws.Cells.Item(301,16).Value = 3400
$ws.Cells.Item(301,17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(301,18).Value = "Provincia de Limarí"
$ws.Cells.Item(301,19).Value = 212
$ws.Cells.Item(301,20).Value = 16

# Row 302: new "1a/2a/3a amarillo" price entry for fecha 44468
$ws.Cells.Item(302,1).Value = 2
$ws.Cells.Item(302,2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(302,3).Value = "Coquimbo"
$ws.Cells.Item(302,4).Value = 44468
$ws.Cells.Item(302,5).Value = 4
$ws.Cells.Item(302,6).Value = "Fruta"
$ws.Cells.Item(302,7).Value = 100102
$ws.Cells.Item(302,8).Value = "Cítricos"
$ws.Cells.Item(302,9).Value = 100102003
$ws.Cells.Item(302,10).Value = "Limón"
$ws.Cells.Item(302,11).Value = "Sin especificar"
$ws.Cells.Item(302,12).Value = "3a amarillo"
$ws.Cells.Item(302,13).Value = 420
$ws.Cells.Item(302,14).Value = 2300
$ws.Cells.Item(302,15).Value = 2500
$ws.Cells.Item(302,16).Value = 2393
$ws.Cells.Item(302,17).Value = "`$/malla 16 kilos"
$ws.Cells.Item(302,18).Value = "Provincia de Limarí"
$ws.Cells.Item(302,19).Value = 150
$ws.Cells.Item(302,20).Value = 16

